$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$cell = $ws.Range("B5")
$old = $cell.Value2
$newText = $old.Replace("= 'Not Reported'`nORDER BY", "= 'Not Reported' AND trt.treatment_id IS NOT NULL`nORDER BY")
$cell.Value2 = $newText
$cell.Font.ThemeColor = 1
Write-Host "DONE"
